$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) labels
$ws.Range("B1").Value = "Email_id"
$ws.Range("D1").Value = "First_Name"
$ws.Range("E1").Value = "Last_Name"

# Update Password column (C2:C6) from numeric placeholders to lowercase username text
$ws.Range("C2").Value = "jay"
$ws.Range("C3").Value = "khan"
$ws.Range("C4").Value = "kovid"
$ws.Range("C5").Value = "simon"
$ws.Range("C6").Value = "tim"

# Update the active selection to match the saved view state
$ws.Range("G12").Select()
